$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above the "مناديل جيب مبلله" row (row 10), shifting
#    that row (and the summary/footer rows below it) down by one.
$ws.Rows.Item(10).Insert()

# 2. Copy the row formatting (now sitting on row 11, the shifted item row)
#    onto the freshly inserted blank row 10 so it matches the other item rows.
$ws.Range("A11:Q11").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$ws.Rows.Item(10).RowHeight = 24.75

# 3. Recreate the merged cells for the new row 10 (lost on insert).
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# 4. Fill in the new item's data: "اختبار حمل بيبي تشك" (pregnancy test).
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "اختبار حمل بيبي تشك "
$ws.Range("H10").Value = "5:0"
$ws.Range("Q10").Value = "1:0"

# These three look like plain numbers ("0", "25.00", "25.0000") - without
# forcing a text format first Excel would silently store them as numeric
# values instead of the text strings the source file uses, so flip the
# cell to text, assign, then restore the original numeric display format
# (copied from the matching column of the row below) so the cell keeps its
# original style.
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "0"
$ws.Range("N10").NumberFormat = "@"
$ws.Range("N10").Value = "25.00"
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "25.0000"

$ws.Range("L11").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("N11").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("P11").Copy()
$ws.Range("P10").PasteSpecial(-4122)

# 5. The "مناديل جيب مبلله" row is now row 11 - bump its item number to 5.
$ws.Range("A11").Value = 5

# 6. Update the page-total row (now row 12) to include the new item's price.
$ws.Range("P12").Value = 212.14

# 7. Update the footer timestamp (now row 13).
$ws.Range("A13").Value = "Sunday, 17 August, 2025 9:45 AM"

Write-Output "edit applied"
